$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel COM constants (not predefined in this host, so spelled out explicitly)
$xlPasteFormats = -4122
$xlRight = -4152

# ---------------------------------------------------------------------------
# 1. New column L = year 2022, mirroring the existing 2014..2021 columns.
# ---------------------------------------------------------------------------

# Header (row 3): same style as the other year headers (e.g. K3).
$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L3").Value = 2022

# Row 4 - "Generation of hazardous waste per person" (computed ratio).
$ws.Range("J4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L4").Formula = "=L5/L6*1000"

# Row 5 - "Hazardous waste generation, thousand tons".
$ws.Range("K5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L5").Value = 12673.2
$ws.Range("L5").HorizontalAlignment = $xlRight

# Row 6 - "Resident population, thousand people".
$ws.Range("K6").Copy() | Out-Null
$ws.Range("L6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L6").Value = 7037.6
$ws.Range("L6").HorizontalAlignment = $xlRight

# Row 7 - "Percentage of neutralized hazardous waste, percent".
$ws.Range("K7").Copy() | Out-Null
$ws.Range("L7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L7").Value = 51.3
$ws.Range("L7").HorizontalAlignment = $xlRight

# Row 8 - "The share of hazardous waste buried, percent".
$ws.Range("K8").Copy() | Out-Null
$ws.Range("L8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L8").Value = 0.1

# ---------------------------------------------------------------------------
# 2. Missing 2014 data points for the two bottom indicators become "-"
#    placeholders, right-aligned like the numeric cells around them.
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "-"
$ws.Range("D7").HorizontalAlignment = $xlRight

$ws.Range("D8").Value = "-"
$ws.Range("D8").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------------
# 3. Selection cursor, matching the saved view of the workbook.
# ---------------------------------------------------------------------------
$ws.Range("N5").Select() | Out-Null
